$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 = "time_taken", styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# time_taken values for data rows 2..102 (one per gene row)
$timeTaken = @("2021-10-05 10:52:16.034014","2021-10-05 10:52:16.034025","2021-10-05 10:52:16.034029","2021-10-05 10:52:16.034031","2021-10-05 10:52:16.034034","2021-10-05 10:52:16.034037","2021-10-05 10:52:16.034040","2021-10-05 10:52:16.034042","2021-10-05 10:52:16.034045","2021-10-05 10:52:16.034047","2021-10-05 10:52:16.034050","2021-10-05 10:52:16.034052","2021-10-05 10:52:16.034055","2021-10-05 10:52:16.034057","2021-10-05 10:52:16.034060","2021-10-05 10:52:16.034062","2021-10-05 10:52:16.034065","2021-10-05 10:52:16.034067","2021-10-05 10:52:16.034070","2021-10-05 10:52:16.034072","2021-10-05 10:52:16.034075","2021-10-05 10:52:16.034077","2021-10-05 10:52:16.034080","2021-10-05 10:52:16.034082","2021-10-05 10:52:16.034085","2021-10-05 10:52:16.034088","2021-10-05 10:52:16.034090","2021-10-05 10:52:16.034093","2021-10-05 10:52:16.034095","2021-10-05 10:52:16.034098","2021-10-05 10:52:16.034100","2021-10-05 10:52:16.034103","2021-10-05 10:52:16.034106","2021-10-05 10:52:16.034108","2021-10-05 10:52:16.034111","2021-10-05 10:52:16.034113","2021-10-05 10:52:16.034116","2021-10-05 10:52:16.034118","2021-10-05 10:52:16.034121","2021-10-05 10:52:16.034123","2021-10-05 10:52:16.034126","2021-10-05 10:52:16.034129","2021-10-05 10:52:16.034131","2021-10-05 10:52:16.034134","2021-10-05 10:52:16.034136","2021-10-05 10:52:16.034139","2021-10-05 10:52:16.034142","2021-10-05 10:52:16.034144","2021-10-05 10:52:16.034147","2021-10-05 10:52:16.034149","2021-10-05 10:52:16.034151","2021-10-05 10:52:16.034154","2021-10-05 10:52:16.034157","2021-10-05 10:52:16.034160","2021-10-05 10:52:16.034162","2021-10-05 10:52:16.034165","2021-10-05 10:52:16.034167","2021-10-05 10:52:16.034170","2021-10-05 10:52:16.034172","2021-10-05 10:52:16.034174","2021-10-05 10:52:16.034177","2021-10-05 10:52:16.034179","2021-10-05 10:52:16.034182","2021-10-05 10:52:16.034185","2021-10-05 10:52:16.034188","2021-10-05 10:52:16.034191","2021-10-05 10:52:16.034194","2021-10-05 10:52:16.034196","2021-10-05 10:52:16.034199","2021-10-05 10:52:16.034201","2021-10-05 10:52:16.034204","2021-10-05 10:52:16.034206","2021-10-05 10:52:16.034209","2021-10-05 10:52:16.034211","2021-10-05 10:52:16.034214","2021-10-05 10:52:16.034216","2021-10-05 10:52:16.034220","2021-10-05 10:52:16.034223","2021-10-05 10:52:16.034226","2021-10-05 10:52:16.034228","2021-10-05 10:52:16.034231","2021-10-05 10:52:16.034234","2021-10-05 10:52:16.034236","2021-10-05 10:52:16.034239","2021-10-05 10:52:16.034241","2021-10-05 10:52:16.034244","2021-10-05 10:52:16.034246","2021-10-05 10:52:16.034249","2021-10-05 10:52:16.034251","2021-10-05 10:52:16.034254","2021-10-05 10:52:16.034256","2021-10-05 10:52:16.034259","2021-10-05 10:52:16.034263","2021-10-05 10:52:16.034266","2021-10-05 10:52:16.034269","2021-10-05 10:52:16.034271","2021-10-05 10:52:16.034274","2021-10-05 10:52:16.034276","2021-10-05 10:52:16.034279","2021-10-05 10:52:16.034282","2021-10-05 10:52:16.034285")

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
